$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.205.13"
$ws.Range("E2").Value = "  -3.40%  "

$ws.Range("D3").Value = "1.614.36"
$ws.Range("E3").Value = "  -2.59%  "

$ws.Range("D4").Value = "'0.9984"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'0.9979"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").Value = "'302.96"
$ws.Range("E6").Value = "  -2.26%  "

$ws.Range("D7").Value = "'0.3783"
$ws.Range("E7").Value = "  -3.26%  "

$ws.Range("D8").Value = "'0.3691"
$ws.Range("E8").Value = "  -4.10%  "

$ws.Range("D9").Value = "'49.12"
$ws.Range("E9").Value = "  -4.24%  "

$ws.Range("D10").Value = "'0.9981"
$ws.Range("E10").Value = "  -0.05%  "

$ws.Range("D11").Value = "'1.285"
$ws.Range("E11").Value = "  -5.69%  "

$ws.Range("D12").Value = "'0.08118"
$ws.Range("E12").Value = "  -4.14%  "

$ws.Range("D13").Value = "'23.27"
$ws.Range("E13").Value = "  -3.61%  "

$ws.Range("D14").Value = "'6.665"
$ws.Range("E14").Value = "  -6.69%  "

$ws.Range("D15").Value = "'7.683"
$ws.Range("E15").Value = "  -2.89%  "

$ws.Range("D16").Value = "'0.00001279"
$ws.Range("E16").Value = "  -3.11%  "

$ws.Range("D17").Value = "1.606.07"
$ws.Range("E17").Value = "  -2.84%  "

$ws.Range("D18").Value = "'91.79"
$ws.Range("E18").Value = "  -3.08%  "

$ws.Range("D19").Value = "'0.06801"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("D20").Value = "'18.52"
$ws.Range("E20").Value = "  -6.60%  "

$ws.Range("D21").Value = "'6.621"
$ws.Range("E21").Value = "  -4.51%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "'13.12"
$ws.Range("E23").Value = "  -3.92%  "

$ws.Range("D24").Value = "23.217.47"
$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("D25").Value = "'2.356"
$ws.Range("E25").Value = "  -5.22%  "

$ws.Range("D26").Value = "'2.940"
$ws.Range("E26").Value = "  -2.67%  "

$ws.Range("E27").Value = "  -4.41%  "

$ws.Range("D28").Value = "'150.99"
$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("D29").Value = "'5.288"
$ws.Range("E29").Value = "  -3.11%  "

$ws.Range("D30").Value = "'132.79"

$ws.Range("D31").Value = "'2.429"
$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").Value = "'7.048"
$ws.Range("E32").Value = "  -10.73%  "

$ws.Range("D33").Value = "1.783.82"
$ws.Range("E33").Value = "  -2.77%  "

$ws.Range("D34").Value = "'1.001"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("D35").Value = "'0.07780"
$ws.Range("E35").Value = "  -4.19%  "

$ws.Range("D36").Value = "'0.02803"
$ws.Range("E36").Value = "  -6.15%  "

$ws.Range("D37").Value = "'6.363"
$ws.Range("E37").Value = "  -6.06%  "

$ws.Range("D38").Value = "'0.2569"
$ws.Range("E38").Value = "  -4.34%  "

$ws.Range("D39").Value = "'10.19"
$ws.Range("E39").Value = "  -6.54%  "

$ws.Range("D40").Value = "'0.08896"

$ws.Range("D41").Value = "'1.402"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("D42").Value = "'0.7228"
$ws.Range("E42").Value = "  -4.71%  "

$ws.Range("D43").Value = "'12.88"
$ws.Range("E43").Value = "  -4.52%  "

$ws.Range("D44").Value = "'16.04"
$ws.Range("E44").Value = "  -2.24%  "

$ws.Range("D45").Value = "'0.6669"
$ws.Range("E45").Value = "  -4.23%  "

$ws.Range("D46").Value = "'2.322"
$ws.Range("E46").Value = "  -5.77%  "

$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").Value = "'3.981"
$ws.Range("E48").Value = "  -2.55%  "

$ws.Range("D49").Value = "'0.08049"
$ws.Range("E49").Value = "  -3.04%  "

$ws.Range("D50").Value = "'131.83"
$ws.Range("E50").Value = "  -2.18%  "

$ws.Range("D51").Value = "'1.176"
$ws.Range("E51").Value = "  -3.74%  "
